$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Anpep"
$ws.Range("C2").Value = "Sele"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.829300666666667
$ws.Range("H2").Value = 14.487902
$ws.Range("I2").Value = 0.02007571491808102
$ws.Range("J2").Value = 0.02007571491808102
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.699506666666667
$ws.Range("N2").Value = 14.09852
$ws.Range("O2").Value = 0.9660495246229048
$ws.Range("P2").Value = 0.9660495246229047
$ws.Range("Q2").Value = 22.69533067833778
$ws.Range("R2").Value = 204.25797610504
$ws.Range("S2").Value = 0.01939413485307712
$ws.Range("T2").Value = 0.01939413485307712

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Anpep"
$ws.Range("C3").Value = "Sele"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.829300666666667
$ws.Range("H3").Value = 14.487902
$ws.Range("I3").Value = 0.02007571491808102
$ws.Range("J3").Value = 0.02007571491808102
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.1651576666666667
$ws.Range("N3").Value = 0.495473
$ws.Range("O3").Value = 0.03395047537709522
$ws.Range("P3").Value = 0.03395047537709522
$ws.Range("Q3").Value = 0.7975960297384445
$ws.Range("R3").Value = 7.178364267646
$ws.Range("S3").Value = 0.0006815800650038927
$ws.Range("T3").Value = 0.0006815800650038927

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Anpep"
$ws.Range("C4").Value = "Sele"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 143.66272
$ws.Range("H4").Value = 430.98816
$ws.Range("I4").Value = 0.597215209850832
$ws.Range("J4").Value = 0.5972152098508319
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.699506666666667
$ws.Range("N4").Value = 14.09852
$ws.Range("O4").Value = 0.9660495246229048
$ws.Range("P4").Value = 0.9660495246229047
$ws.Range("Q4").Value = 675.1439103914668
$ws.Range("R4").Value = 6076.2951935232
$ws.Range("S4").Value = 0.5769394695739646
$ws.Range("T4").Value = 0.5769394695739645

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Anpep"
$ws.Range("C5").Value = "Sele"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 143.66272
$ws.Range("H5").Value = 430.98816
$ws.Range("I5").Value = 0.597215209850832
$ws.Range("J5").Value = 0.5972152098508319
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1651576666666667
$ws.Range("N5").Value = 0.495473
$ws.Range("O5").Value = 0.03395047537709522
$ws.Range("P5").Value = 0.03395047537709522
$ws.Range("Q5").Value = 23.72699962218667
$ws.Range("R5").Value = 213.54299659968
$ws.Range("S5").Value = 0.02027574027686743
$ws.Range("T5").Value = 0.02027574027686742

$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Anpep"
$ws.Range("C6").Value = "Sele"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 84.73579966666667
$ws.Range("H6").Value = 254.207399
$ws.Range("I6").Value = 0.3522521944440867
$ws.Range("J6").Value = 0.3522521944440867
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.699506666666667
$ws.Range("N6").Value = 14.09852
$ws.Range("O6").Value = 0.9660495246229048
$ws.Range("P6").Value = 0.9660495246229047
$ws.Range("Q6").Value = 398.2164554388311
$ws.Range("R6").Value = 3583.94809894948
$ws.Range("S6").Value = 0.340293064990085
$ws.Range("T6").Value = 0.340293064990085

$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Anpep"
$ws.Range("C7").Value = "Sele"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 84.73579966666667
$ws.Range("H7").Value = 254.207399
$ws.Range("I7").Value = 0.3522521944440867
$ws.Range("J7").Value = 0.3522521944440867
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.1651576666666667
$ws.Range("N7").Value = 0.495473
$ws.Range("O7").Value = 0.03395047537709522
$ws.Range("P7").Value = 0.03395047537709522
$ws.Range("Q7").Value = 13.99476695608078
$ws.Range("R7").Value = 125.952902604727
$ws.Range("S7").Value = 0.01195912945400172
$ws.Range("T7").Value = 0.01195912945400172

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Anpep"
$ws.Range("C8").Value = "Sele"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 7.326535333333333
$ws.Range("H8").Value = 21.979606
$ws.Range("I8").Value = 0.03045688078700028
$ws.Range("J8").Value = 0.03045688078700028
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.699506666666667
$ws.Range("N8").Value = 14.09852
$ws.Range("O8").Value = 0.9660495246229048
$ws.Range("P8").Value = 0.9660495246229047
$ws.Range("Q8").Value = 34.43110164256889
$ws.Range("R8").Value = 309.87991478312
$ws.Range("S8").Value = 0.02942285520577811
$ws.Range("T8").Value = 0.0294228552057781

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Anpep"
$ws.Range("C9").Value = "Sele"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 7.326535333333333
$ws.Range("H9").Value = 21.979606
$ws.Range("I9").Value = 0.03045688078700028
$ws.Range("J9").Value = 0.03045688078700028
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1651576666666667
$ws.Range("N9").Value = 0.495473
$ws.Range("O9").Value = 0.03395047537709522
$ws.Range("P9").Value = 0.03395047537709522
$ws.Range("Q9").Value = 1.210033480404222
$ws.Range("R9").Value = 10.890301323638
$ws.Range("S9").Value = 0.001034025581222177
$ws.Range("T9").Value = 0.001034025581222177
